$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = "State of the art in image manipulation (stylegan)`nhttps://www.catalyzex.com/pape"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2708"

$ws.Range("D29").Value = "프로메디우스"

$ws.Range("D32").Value = "Gumbel Softmax : 카테고리형 변수 변환"
